$wb = $excel.ActiveWorkbook

# The handoff report text moved on from "Ready for handoff" to "In Translation".
# Update every sheet's Status column(s) and re-fit those columns to the new
# (shorter) text, matching how the report's refresh/archive step narrows them.

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$fitWidth  = 12.5

# Overview sheet: per-locale status columns "zh-cn" (E) and "de-de" (F)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.UsedRange.Replace($oldStatus, $newStatus) | Out-Null
$wsOverview.Columns.Item(5).ColumnWidth = $fitWidth
$wsOverview.Columns.Item(6).ColumnWidth = $fitWidth

# zh-cn sheet: "Status" column (C)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.UsedRange.Replace($oldStatus, $newStatus) | Out-Null
$wsZhCn.Columns.Item(3).ColumnWidth = $fitWidth

# de-de sheet: "Status" column (C)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.UsedRange.Replace($oldStatus, $newStatus) | Out-Null
$wsDeDe.Columns.Item(3).ColumnWidth = $fitWidth
